# recipes_total.xlsx -- "aggiunta ricette e sistemazione pt.1"
# Replace the 5-recipe placeholder table with the final 8-recipe table,
# shrink the sheet from 14 to 11 rows, fix up number formats / wrap text
# on the re-shuffled rows, retint the title font, resize the columns and
# restore the saved zoom / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- number-format strings that must resolve to the workbook's existing
#     custom numFmts (164 / 165 / 6 / 8) instead of minting new ones ---
$fmt164 = '#,##0.00\ "€"'
$fmt165 = '#,##0.0\ "€"'
$fmt6   = '#,##0\ "€";[Red]\-#,##0\ "€"'
$fmt8   = '#,##0.00\ "€";[Red]\-#,##0.00\ "€"'

# --- shrink the sheet: it now ends at row 11, not row 14 ---
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()

# --- row 1 : Cheesy Egg Stuffed Zucchini Boats ---
# (keep A1's existing "title" style slot -- only its font changes)
$ws.Range("A1").Font.Name = "Aptos Narrow"
$ws.Range("A1").Value = "Cheesy Egg Stuffed Zucchini Boats"
$ws.Range("B1").ClearFormats()
$ws.Range("B1").Value = "zucchini, egg, potato, cheese"
$ws.Range("C1").Value = 30
$ws.Range("D1").ClearFormats()
$ws.Range("D1").Value = "easy"
$ws.Range("E1").ClearFormats()
$ws.Range("E1").Value = "Preheat your oven to 190°C. Place the zucchini halves on a baking sheet. Sprinkle with salt and pepper. In a bowl, mix grated potato, shredded cheese, chopped parsley, salt, and pepper. Fill each zucchini half with the potato and cheese mixture. Carefully crack an egg into each filled zucchini half. Bake in the preheated oven for 15-20 minutes, or until the egg whites are set but the yolks are still runny. Serve hot, garnished with additional chopped parsley if desired."

# --- row 2 : Vegetable frittata ---
$ws.Range("A2").ClearFormats()
$ws.Range("A2").Value = "Vegetable frittata"
$ws.Range("B2").ClearFormats()
$ws.Range("B2").NumberFormat = $fmt164
$ws.Range("B2").Value = "carrot, zucchini, egg"
$ws.Range("C2").Value = 30
$ws.Range("D2").ClearFormats()
$ws.Range("D2").Value = "easy"
$ws.Range("E2").ClearFormats()
$ws.Range("E2").Value = "In a large oven-safe skillet, heat olive oil over medium heat. Add the grated carrot and diced zucchini to the skillet and cook until they start to soften, about 5-7 minutes. In a bowl, beat the eggs and season with salt and pepper. Pour the beaten eggs over the vegetables in the skillet. Let it cook for 2-3 minutes until the frittata is set. Then serve hot."

# --- row 3 : Loaded Potato Skins ---
$ws.Range("A3").ClearFormats()
$ws.Range("A3").Value = "Loaded Potato Skins"
$ws.Range("B3").ClearFormats()
$ws.Range("B3").NumberFormat = $fmt165
$ws.Range("B3").Value = "potato, cheese"
$ws.Range("C3").Value = 60
$ws.Range("D3").ClearFormats()
$ws.Range("D3").Value = "moderate"
$ws.Range("E3").ClearFormats()
$ws.Range("E3").NumberFormat = $fmt165
$ws.Range("E3").Value = "Preheat your oven to 200°C. Scrub the potatoes clean and pierce each potato several times with a fork. Place the potatoes directly on the oven rack and bake until they are tender when pierced with a fork. Remove the potatoes from the oven and once they are cool enough to handle, slice each potato in half lengthwise. Scoop out the flesh, leaving about ¼ inch of potato on the skin. Place the potato skins on a baking sheet and fill each potato skin with shredded cheese. Return the filled potato skins to the oven and bake for an additional 10-15 minutes, or until the cheese is melted and bubbly. Remove from the oven and serve."

# --- row 4 : Potato and Cheese Croquettes ---
$ws.Range("A4").ClearFormats()
$ws.Range("A4").Value = "Potato and Cheese Croquettes"
$ws.Range("B4").ClearFormats()
$ws.Range("B4").NumberFormat = $fmt6
$ws.Range("B4").Value = "potato, cheese, egg"
$ws.Range("C4").Value = 45
$ws.Range("D4").ClearFormats()
$ws.Range("D4").Value = "moderate"
$ws.Range("E4").ClearFormats()
$ws.Range("E4").Value = "Mash the boiled potatoes in a bowl. Mix in the shredded cheese, salt, and pepper. Take small portions of the potato mixture and shape them into croquettes. Dip each croquette into the beaten eggs, then coat evenly with breadcrumbs. Heat oil in a frying pan over medium heat. Fry the croquettes until golden brown and crispy on all sides. Remove from the oil and place on a paper towel-lined plate to drain excess oil. Serve hot with your favorite dipping sauce."

# --- row 5 : Zucchini Cheese Rollbacks ---
$ws.Range("A5").ClearFormats()
$ws.Range("A5").Value = "Zucchini Cheese Rollbacks"
$ws.Range("B5").ClearFormats()
$ws.Range("B5").NumberFormat = $fmt6
$ws.Range("B5").Value = "zucchini, cheese"
$ws.Range("C5").Value = 45
$ws.Range("D5").ClearFormats()
$ws.Range("D5").Value = "moderate"
$ws.Range("E5").ClearFormats()
$ws.Range("E5").Value = "Preheat your oven to 190°C. Grease a baking dish with olive oil and set aside. Trim the ends of the zucchinis and slice them lengthwise into thin strips, about 1/4 inch thick. In a bowl, combine the shredded cheese, grated Parmesan cheese, garlic powder, salt, and pepper. Mix until well combined. Lay out the zucchini slices on a flat surface. Spread the cheese mixture evenly over each zucchini slice. Carefully roll up each zucchini slice with the cheese mixture inside to create rollbacks. Place the zucchini rollbacks seam side down in the prepared baking dish. Bake in the preheated oven until the zucchini is tender and the cheese is melted and bubbly. Remove from the oven and let cool slightly before serving."

# --- row 6 : Potato and Carrot Balls ---
$ws.Range("A6").ClearFormats()
$ws.Range("A6").Value = "Potato and Carrot Balls"
$ws.Range("B6").ClearFormats()
$ws.Range("B6").NumberFormat = $fmt6
$ws.Range("B6").Value = "potato, carrot"
$ws.Range("C6").Value = 20
$ws.Range("D6").ClearFormats()
$ws.Range("D6").Value = "easy"
$ws.Range("E6").ClearFormats()
$ws.Range("E6").Value = "Boil the potatoes, scrub them with a fork and add grated carrot. Add flour, salt, paprik and pepper to the mixture. Take small portions of the mixture and shape them into little balls.  Heat oil in a frying pan over medium heat. Fry the balls until golden brown and crispy on all sides. Remove from the oil and place on a paper towel-lined plate to drain excess oil. "

# --- row 7 : Carrot chips ---
$ws.Range("A7").ClearFormats()
$ws.Range("A7").Value = "Carrot chips"
$ws.Range("B7").ClearFormats()
$ws.Range("B7").NumberFormat = $fmt6
$ws.Range("B7").Value = "carrot"
$ws.Range("C7").Value = 30
$ws.Range("D7").ClearFormats()
$ws.Range("D7").Value = "easy"
$ws.Range("E7").ClearFormats()
$ws.Range("E7").Value = "Cut the carrots into circles, wash them under cold water and dry them. Add salt, oil and paprik to the carrots and mix it all. Place the carrots on the oven rack without overlapping them and then cook untile they aren't crispy and dry. Add some more salt and serve with a yogurt sauce."

# --- row 8 : Zucchini Pesto ---
$ws.Range("A8").ClearFormats()
$ws.Range("A8").WrapText = $true
$ws.Range("A8").Value = "Zucchini Pesto"
$ws.Range("B8").ClearFormats()
$ws.Range("B8").NumberFormat = $fmt6
$ws.Range("B8").Value = "zucchini, cheese"
$ws.Range("C8").Value = 10
$ws.Range("D8").ClearFormats()
$ws.Range("D8").Value = "easy"
$ws.Range("E8").ClearFormats()
$ws.Range("E8").Value = "Wash zucchini, trim the ends and cut them into pieces. Transfer the pieces into a mixer and add olice oli, salt, parmisan cheese and basil. Mix the ingredients and serve."
$ws.Rows.Item(8).RowHeight = 17

# --- rows 9-11 : formatting-only placeholder cells, no content ---
$ws.Range("B9").ClearContents()
$ws.Range("B9").NumberFormat = $fmt8
$ws.Range("B10").ClearContents()
$ws.Range("B10").NumberFormat = $fmt8
$ws.Range("A11").ClearContents()
$ws.Range("A11").ClearFormats()
$ws.Range("B11").ClearContents()
$ws.Range("B11").NumberFormat = $fmt6

# --- column widths ---
$ws.Columns.Item(1).ColumnWidth = 30.666666666666668
$ws.Columns.Item(2).ColumnWidth = 23.166666666666668
$ws.Columns.Item(3).ColumnWidth = 8.666666666666666
$ws.Columns.Item(4).ColumnWidth = 11.666666666666666
$ws.Columns.Item(5).ColumnWidth = 35.333333333333336

# --- view: zoom + selection ---
$excel.ActiveWindow.Zoom = 162
$ws.Range("B1").Select()
